$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value without Excel re-interpreting
# numeric-looking strings (e.g. "20.70") as numbers, which would silently
# drop significant trailing zeros. We flip the cell to Text format just
# long enough to assign the literal, then restore the original (General,
# un-styled) formatting by pasting formats in from an always-blank donor
# cell, so no stray style index is left behind on the cell.
function Set-TextValue($cell, $text) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $text
    $ws.Range("ZZ1").Copy() | Out-Null
    $ws.Range($cell).PasteSpecial(-4122) | Out-Null
}

# Row 2
$ws.Range("D2").Value = "37.413.93"
$ws.Range("E2").Value = "  -1.03%  "

# Row 3
$ws.Range("D3").Value = "2.050.51"
$ws.Range("E3").Value = "  -1.50%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
Set-TextValue "D5" "229.25"
$ws.Range("E5").Value = "  -1.72%  "

# Row 6
$ws.Range("E6").Value = "  -1.36%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
Set-TextValue "D8" "57.17"
$ws.Range("E8").Value = "  -2.44%  "

# Row 9
Set-TextValue "D9" "0.386"
$ws.Range("E9").Value = "  -1.71%  "

# Row 10
Set-TextValue "D10" "0.0791"
$ws.Range("E10").Value = "  +0.61%  "

# Row 11
$ws.Range("E11").Value = "  -1.90%  "

# Row 12
Set-TextValue "D12" "14.77"
$ws.Range("E12").Value = "  -0.88%  "

# Row 13
$ws.Range("D13").Value = "2.351.33"
$ws.Range("E13").Value = "  -1.49%  "

# Row 14
Set-TextValue "D14" "20.70"
$ws.Range("E14").Value = "  -2.07%  "

# Row 15
$ws.Range("E15").Value = "  -2.99%  "

# Row 16
$ws.Range("E16").Value = "  -0.57%  "

# Row 17
$ws.Range("D17").Value = "2.047.22"
$ws.Range("E17").Value = "  -1.90%  "

# Row 18
$ws.Range("D18").Value = "37.286.84"
$ws.Range("E18").Value = "  -1.23%  "

# Row 19
Set-TextValue "D19" "6.10"

# Row 20
Set-TextValue "D20" "69.63"
$ws.Range("E20").Value = "  -2.60%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0830"
$ws.Range("E21").Value = "  -1.36%  "

# Row 22
Set-TextValue "D22" "226.47"
$ws.Range("E22").Value = "  -1.39%  "

# Row 23
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
Set-TextValue "D25" "2.31"
$ws.Range("E25").Value = "  -3.72%  "

# Row 26
$ws.Range("E26").Value = "  -1.00%  "

# Row 27
Set-TextValue "D27" "168.41"
$ws.Range("E27").Value = "  -2.10%  "

# Row 28
$ws.Range("E28").Value = "  -6.05%  "

# Row 29
Set-TextValue "D29" "18.99"
$ws.Range("E29").Value = "  -2.50%  "

# Row 30
$ws.Range("E30").Value = "  -3.27%  "

# Row 31
$ws.Range("E31").Value = "  -2.36%  "

# Row 32
$ws.Range("E32").Value = "  -3.94%  "

# Row 33
$ws.Range("E33").Value = "  -2.70%  "

# Row 34
Set-TextValue "D34" "4.58"
$ws.Range("E34").Value = "  -2.07%  "

# Row 35
$ws.Range("E35").Value = "  -1.29%  "

# Row 36
$ws.Range("E36").Value = "  +1.51%  "

# Row 37
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("E38").Value = "  -4.93%  "

# Row 39
Set-TextValue "D39" "5.33"
$ws.Range("E39").Value = "  -1.77%  "

# Row 40
Set-TextValue "D40" "0.0223"
$ws.Range("E40").Value = "  -5.12%  "

# Row 41
Set-TextValue "D41" "17.23"
$ws.Range("E41").Value = "  +1.46%  "

# Row 42
$ws.Range("E42").Value = "  -1.19%  "

# Row 43
$ws.Range("D43").Value = "1.474.53"
$ws.Range("E43").Value = "  +2.02%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "96.44"
$ws.Range("E44").Value = "  -4.71%  "

# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D45" "0.0940"
$ws.Range("E45").Value = "  -3.33%  "

# Row 46
$ws.Range("E46").Value = "  +0.70%  "

# Row 47
Set-TextValue "D47" "1.03"
$ws.Range("E47").Value = "  -4.13%  "

# Row 48
Set-TextValue "D48" "3.93"
$ws.Range("E48").Value = "  -4.42%  "

# Row 49
$ws.Range("E49").Value = "  -2.52%  "

# Row 50
$ws.Range("E50").Value = "  -2.21%  "

# Row 51
$ws.Range("D51").Value = "2.240.36"
$ws.Range("E51").Value = "  -1.41%  "
